# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Anima Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4796.273
$ws.Range("J131").Value = 5752.6665
$ws.Range("L131").Value = 17257.9995
$ws.Range("N131").Value = -27337.9995
$ws.Range("H132").Value = 2532.8086
$ws.Range("I132").Value = 2424.8333
$ws.Range("K132").Value = 7274.499899999999
$ws.Range("M132").Value = -4744.499899999999
$ws.Range("H137").Value = 1284.84
$ws.Range("I137").Value = 1197
$ws.Range("J137").Value = 1416.6
$ws.Range("K137").Value = 3591
$ws.Range("L137").Value = 4249.799999999999
$ws.Range("M137").Value = -1041
$ws.Range("N137").Value = -9349.799999999999
$ws.Range("H138").Value = 1289.35
$ws.Range("I138").Value = 510.83334
$ws.Range("J138").Value = 2007.9807
$ws.Range("K138").Value = 1532.50002
$ws.Range("L138").Value = 6023.9421
$ws.Range("M138").Value = 3607.49998
$ws.Range("N138").Value = -16303.9421

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5802.6
$ws.Range("I26").Value = 4997
$ws.Range("J26").Value = 6004
$ws.Range("K26").Value = 4997
$ws.Range("L26").Value = 6004
$ws.Range("M26").Value = -4667
$ws.Range("N26").Value = -6664
$ws.Range("H32").Value = 940629
$ws.Range("I32").Value = 1065186.8
$ws.Range("K32").Value = 1065186.8
$ws.Range("M32").Value = -1064899.8
$ws.Range("H61").Value = 2320.7795
$ws.Range("I61").Value = 1944.9524
$ws.Range("K61").Value = 1944.9524
$ws.Range("M61").Value = -1732.9524
$ws.Range("H74").Value = 906.39215
$ws.Range("I74").Value = 653.08105
$ws.Range("J74").Value = 1575.8572
$ws.Range("K74").Value = 653.08105
$ws.Range("L74").Value = 1575.8572
$ws.Range("M74").Value = 220.91895
$ws.Range("N74").Value = -3323.8572
$ws.Range("H77").Value = 906.39215
$ws.Range("I77").Value = 653.08105
$ws.Range("J77").Value = 1575.8572
$ws.Range("K77").Value = 3265.40525
$ws.Range("L77").Value = 7879.286
$ws.Range("M77").Value = 1102.59475
$ws.Range("N77").Value = -16615.286
$ws.Range("H136").Value = 2320.7795
$ws.Range("I136").Value = 1944.9524
$ws.Range("K136").Value = 5834.857199999999
$ws.Range("M136").Value = -3284.857199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = ""
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = ""
$ws.Range("N79").Value = 0
$ws.Range("H112").Value = 98469
$ws.Range("J112").Value = 98469
$ws.Range("L112").Value = 98469
$ws.Range("N112").Value = -101423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 234.66667
$ws.Range("I7").Value = 286
$ws.Range("J7").Value = 170.5
$ws.Range("K7").Value = 286
$ws.Range("L7").Value = 170.5
$ws.Range("M7").Value = -173
$ws.Range("N7").Value = -396.5
$ws.Range("H31").Value = 4841.9556
$ws.Range("I31").Value = 1391.7894
$ws.Range("J31").Value = 7363.231
$ws.Range("K31").Value = 1391.7894
$ws.Range("L31").Value = 7363.231
$ws.Range("M31").Value = -1096.7894
$ws.Range("N31").Value = -7953.231
$ws.Range("H34").Value = 4841.9556
$ws.Range("I34").Value = 1391.7894
$ws.Range("J34").Value = 7363.231
$ws.Range("K34").Value = 1391.7894
$ws.Range("L34").Value = 7363.231
$ws.Range("M34").Value = -1189.7894
$ws.Range("N34").Value = -7767.231
$ws.Range("H58").Value = 1114.6046
$ws.Range("I58").Value = 890.5909
$ws.Range("J58").Value = 1349.2858
$ws.Range("K58").Value = 890.5909
$ws.Range("L58").Value = 1349.2858
$ws.Range("M58").Value = -687.5909
$ws.Range("N58").Value = -1755.2858
$ws.Range("H134").Value = 3377.0613
$ws.Range("I134").Value = 3686.8
$ws.Range("J134").Value = 2000.4445
$ws.Range("K134").Value = 11060.4
$ws.Range("L134").Value = 6001.333500000001
$ws.Range("M134").Value = -8525.400000000001
$ws.Range("N134").Value = -11071.3335
$ws.Range("H136").Value = 1114.6046
$ws.Range("I136").Value = 890.5909
$ws.Range("J136").Value = 1349.2858
$ws.Range("K136").Value = 2671.7727
$ws.Range("L136").Value = 4047.8574
$ws.Range("M136").Value = -121.7727
$ws.Range("N136").Value = -9147.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 4058.7778
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 4058.7778
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = ""
$ws.Range("M102").Value = 12176.3334
$ws.Range("N102").Value = -17044.3334
$ws.Range("H113").Value = 669.74286
$ws.Range("I113").Value = 549.3
$ws.Range("J113").Value = 830.3333
$ws.Range("K113").Value = 1647.9
$ws.Range("L113").Value = 2490.9999
$ws.Range("M113").Value = 522.1000000000001
$ws.Range("N113").Value = -6830.9999
$ws.Range("H131").Value = 2720.348
$ws.Range("J131").Value = 3008.7378
$ws.Range("L131").Value = 9026.213400000001
$ws.Range("N131").Value = -19106.2134

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1760.9524
$ws.Range("I102").Value = 1735.8422
$ws.Range("K102").Value = 1735.8422
$ws.Range("M102").Value = -113.8422

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 2552506
$ws.Range("J25").Value = 2552506
$ws.Range("L25").Value = 2552506
$ws.Range("N25").Value = -2552966
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21288
$ws.Range("H132").Value = 2609.383
$ws.Range("I132").Value = 2431.394
$ws.Range("J132").Value = 3028.9285
$ws.Range("K132").Value = 7294.181999999999
$ws.Range("L132").Value = 9086.7855
$ws.Range("M132").Value = -4764.181999999999
$ws.Range("N132").Value = -14146.7855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 95000
$ws.Range("J68").Value = 95000
$ws.Range("L68").Value = 95000
$ws.Range("N68").Value = -96622
$ws.Range("H70").Value = 50900
$ws.Range("J70").Value = 50900
$ws.Range("L70").Value = 50900
$ws.Range("N70").Value = -51530
$ws.Range("H71").Value = 95000
$ws.Range("J71").Value = 95000
$ws.Range("L71").Value = 285000
$ws.Range("N71").Value = -293112
$ws.Range("H73").Value = 50900
$ws.Range("J73").Value = 50900
$ws.Range("L73").Value = 50900
$ws.Range("N73").Value = -53084
$ws.Range("H132").Value = 5209944
$ws.Range("I132").Value = 1899.75
$ws.Range("K132").Value = 5699.25
$ws.Range("M132").Value = -3169.25
